# Apply trade #86 close: 2026-02-17 09:03:01 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.71
$summary.Range("B4").Value = -0.28
$summary.Range("B5").Value = -0.07000000000000001
$summary.Range("B6").Value = 86
$summary.Range("B8").Value = 34
$summary.Range("B9").Value = 43.02

# ---------------------------------------------------------------------------
# Strategy Status sheet
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.70999999999999
$status.Range("D4").Value = 86
$status.Range("E4").Value = -0.28
$status.Range("F4").Value = -0.29
$status.Range("G4").Value = 43.02

# ---------------------------------------------------------------------------
# Helper to append the new trade row (#86) to a trade-log sheet
# ---------------------------------------------------------------------------
function Add-TradeRow87($ws) {
    $ws.Range("A87").Value = 86
    $ws.Range("C87").Value = "09:02:55"
    $ws.Range("D87").Value = "MarketMaking"
    $ws.Range("E87").Value = "UP"
    $ws.Range("F87").Value = 0.67
    $ws.Range("G87").Value = 0.489227
    $ws.Range("H87").Value = "CLOSED"
    $ws.Range("I87").Value = -26.981
    $ws.Range("J87").Value = -0.18
    $ws.Range("K87").Value = 99.70999999999999
    $ws.Range("L87").Value = 0
    $ws.Range("M87").Value = 0
    $ws.Range("N87").Value = 0.6
    $ws.Range("O87").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P87").Value = "early_exit"
    $ws.Range("Q87").Value = 0.14

    # "2026-02-17" looks like a date, and assigning it straight to .Value
    # makes Excel auto-convert it to a date serial number + date style.
    # Route it through a text formula and paste the computed value back in
    # so the cell ends up holding the literal string with no extra style.
    $ws.Range("B87").Formula = '="2026-02-17"'
    $ws.Range("B87").Copy() | Out-Null
    $ws.Range("B87").PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow87 $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow87 $marketMaking
